# Rewrite the player roster table (columns A-C, rows 2-19) with updated
# player / position / team data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ty Jerome"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Cleveland Cavaliers"

$ws.Range("A3").Value = "Duncan Robinson"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Miami Heat"

$ws.Range("A4").Value = "Quentin Grimes"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "Dallas Mavericks"

$ws.Range("A5").Value = "Harrison Barnes"
$ws.Range("B5").Value = "SF,PF"
$ws.Range("C5").Value = "San Antonio Spurs"

$ws.Range("A6").Value = "Brice Sensabaugh"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Utah Jazz"

$ws.Range("A7").Value = "Cameron Johnson"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Brooklyn Nets"

$ws.Range("A8").Value = "Anthony Davis"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Los Angeles Lakers"

$ws.Range("A9").Value = "Bam Adebayo"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Miami Heat"

$ws.Range("A10").Value = "Isaiah Hartenstein"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Oklahoma City Thunder"

$ws.Range("A11").Value = "Julius Randle"
$ws.Range("B11").Value = "PF"
$ws.Range("C11").Value = "Minnesota Timberwolves"

$ws.Range("A12").Value = "Jared McCain"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Philadelphia 76ers"

$ws.Range("A13").Value = "Cade Cunningham"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Detroit Pistons"

$ws.Range("A14").Value = "Damian Lillard"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Milwaukee Bucks"

$ws.Range("A15").Value = "Brandon Miller"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Charlotte Hornets"

$ws.Range("A16").Value = "Bilal Coulibaly"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Washington Wizards"

$ws.Range("A17").Value = "LaMelo Ball"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Charlotte Hornets"

$ws.Range("A18").Value = "Brandon Ingram"
$ws.Range("B18").Value = "SG,SF,PF"
$ws.Range("C18").Value = "New Orleans Pelicans"

$ws.Range("A19").Value = "Derrick White"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Boston Celtics"
